# "Atualização Diária da Ata de reuniões"
# Fills in the new Sprint meeting-minutes row (row 14, day 11/10) that was
# added to the "Atas de reuniões" sheet, reusing the same visual layout
# (borders/fill/font/number format) as the previous entry in row 12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the cell formatting from row 12 (same column layout: Dia / Início /
# Término / Membros Presentes / Membros Ausentes / Pautas Discutidas /
# Assuntos Decididos) onto row 14 so the new row gets matching borders,
# fill, fonts and number formats.
$ws.Range("A12:G12").Copy()
$ws.Range("A14:G14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new meeting data for row 14. Values are entered in the same
# order they first appear so the shared-string table is built in the
# expected sequence.
$ws.Cells.Item(14, 5).Value = "Davi"
$ws.Cells.Item(14, 4).Value = "Arthur`nDaniel`nLeonardo`nMatteus`nPedro"
$ws.Cells.Item(14, 7).Value = "Definimos como sera desenvolvido a modelagem logíca e script do banco de dados (Métricas dos sensores ) e quem desenvolvera;"
$ws.Cells.Item(14, 6).Value = "Compartilhamos experiencia sobre o desenvolvimento do projeto;  definimos como será feito o banco de dados (Métricas dos sensores) ;"
$ws.Cells.Item(14, 2).Value = 0.43055555555555558
$ws.Cells.Item(14, 3).Value = 0.4375

# Resize the row to fit the new content.
$ws.Rows.Item(14).RowHeight = 142.5

# Update the sheet view to reflect the scrolled/zoomed/selected state after
# the edit (window scrolled down to show row 10+, zoomed to 85%, with the
# newly filled-in cell selected).
$win = $ws.Application.ActiveWindow
$win.Zoom = 85
$win.ScrollRow = 10
$win.ScrollColumn = 1
[void]$ws.Range("G14").Select()

Write-Host "Row 14 updated"
